$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9163179916317992
$ws.Range("B3").Value = 0.9180327868852459
$ws.Range("B4").Value = 0.9145299145299145
$ws.Range("B5").Value = 0.9180327868852459
$ws.Range("B6").Value = 0.9145299145299145
$ws.Range("B7").Value = 0.9676334594367381
